# Sprint backlog 1 corrections
# ------------------------------------------------------------------
# 1. Resize the backlog table's grid columns.
# 2. Clean up "able to  send message" (merge runs / drop proofErr marks,
#    text itself is unchanged).
# 3. "Test adding group chat feature" -> "Test sidebar and adding group
#    chat feature".
# 4. Day-6 estimate for that same task row: 1 -> 2.
# 5. Clean up "... to be able to  See what I m texting in real time
#    so that  I can check what I texted before sending the text"
#    (merge runs / drop proofErr marks; a stray space before "so that"
#    is also dropped, matching the author's edit).
# 6. "Write code for hard coded text" -> "Write code to show the name
#    of a sender of the text".
# ------------------------------------------------------------------

$d = $word.ActiveDocument

# 1. Table grid column widths (values are dxa/20 = points).
$t = $d.Tables.Item(1)
$widths = @(504, 3962, 2311, 709, 709, 784, 784, 709, 709, 709)
for ($i = 1; $i -le $widths.Length; $i++) {
    $t.Columns.Item($i).Width = $widths[$i - 1] / 20.0
}

# 2. "able to  send message" - same text, just re-run/merge it so the
#    engine drops the now-redundant proofErr wrappers.
$d.Content.Find.Execute("able to  send message", $true, $false, $false, `
    $false, $false, $true, 1, $false, "able to  send message", 2) | Out-Null

# 3. Insert "sidebar and " into the task description.
$d.Content.Find.Execute("Test adding group chat feature", $true, $false, `
    $false, $false, $false, $true, 1, $false, `
    "Test sidebar and adding group chat feature", 2) | Out-Null

# 4. Bump the Day 6 estimate (column 9) on that task's row from 1 to 2.
#    (direct Range.Text assignment - Find/Replace on a cell-scoped Range
#    can jump outside the cell in this engine, so avoid it here.)
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    $taskCell = $t.Cell($r, 3)
    if ($taskCell.Range.Text.TrimEnd([char]7, [char]13) -eq "Test sidebar and adding group chat feature") {
        $dayCell = $t.Cell($r, 9)
        if ($dayCell.Range.Text.TrimEnd([char]7, [char]13) -eq "1") {
            $dayCell.Range.Text = "2"
        }
    }
}

# 5a. Merge the "to be able ... real time" run and drop the trailing
#     space that used to separate it from the bold "so that " run.
$d.Content.Find.Execute(" to be able to  See what I m texting in real time ", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    " to be able to  See what I m texting in real time", 2) | Out-Null

# 5b. Merge the trailing "I can check..." run (text unchanged).
$d.Content.Find.Execute(" I can check what I texted before sending the text", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    " I can check what I texted before sending the text", 2) | Out-Null

# 6. "Write code for hard coded text" -> "Write code to show the name
#    of a sender of the text".
$d.Content.Find.Execute("Write code for hard coded text", $true, $false, `
    $false, $false, $false, $true, 1, $false, `
    "Write code to show the name of a sender of the text", 2) | Out-Null
